$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Add Chinatown 4K as row 35 (A, C, B, D entry order)
$ws.Range("A35").Value = "Chinatown 4K"
$ws.Range("C35").Value = "https://www.amazon.de/-/en/Jack-Nicholson/dp/B0D9RR7X5P/ref=sr_1_1?crid=3SOGPKQV9UITI&dib=eyJ2IjoiMSJ9.CE3YLu_NPn6syLZ7GBcsz4qEcsjlEm2sGv0Vx-b7F2uh6CyBYY02vTsRT-tQ9-j8GCKpYaRm337tKkSfCi1cgChi-DtEAHPzSOAk0FLkwSUhd-fdWDJTY2raBshE9BFBBSu9syQRVX3GWdFUDzLm0FdQr6ZKHJLQ-cNyE1pQ9ixel8dQmSeThB3EO64eKs-3ly0aaDpa8Zpsv92RYGsW6b3Keao6AC6luXh0j9r_dUk.I3PWPbcxUffDOD-5GzUup630T-8tHzjgm-UGm2IObjA&dib_tag=se&keywords=chinatown+4k&qid=1721912692&sprefix=chinatown+4k%2Caps%2C120&sr=8-1"
$ws.Range("B35").Value = "https://m.media-amazon.com/images/I/71Q2LnsKtWL._SY445_.jpg"
$ws.Range("D35").Value = "21.73 EUR"

# Step 2: Insert a new row at 35, shifting Chinatown down to row 36; fill new row 35 with Abyss 4K (A, B, C, D entry order)
$ws.Rows("35").Insert()
$ws.Range("A35").Value = "Abyss 4K"
$ws.Range("B35").Value = "https://m.media-amazon.com/images/I/810finUC6+L._SX342_.jpg"
$ws.Range("C35").Value = "https://www.amazon.de/-/en/James-Cameron/dp/B0CVKVH3VN/ref=sr_1_1?crid=1EPHQH2T7H52Q&dib=eyJ2IjoiMSJ9.CzkHBhjfMxojXMbBFhQoede9y_v-aB_mhyhLpWIXsRkXeCy__wBLbyq-kjyD_iNQ87ikGF216HZox5M3Vg6jVh5QAEdFBLzrHvVfduidflv_hw19eU7iTmETWD3EAlWHQAYGLQzSm3FfEKV3_uzjc64cbdhBH5bnOVvHd4aMZ5lpHPyLqOgpz3vWu1NNbI_mdterRh-R5jlhcOh2EggHnOfEM4LQl8uH2WVbK8T8TK8.dlkA6Kxh3OzGsrJ2qyOFFQAaGkeGbXDoqseHWzRqQhc&dib_tag=se&keywords=abyss+4k&qid=1721911150&sprefix=abyss+4%2Caps%2C303&sr=8-1"
$ws.Range("D35").Value = "30 EUR"

# Step 3: Add Paprika 4K as row 37 (A, B, C, D entry order)
$ws.Range("A37").Value = "Paprika 4K"
$ws.Range("B37").Value = "https://www.wog.ch/nas/cover_large/4a/4k_paprikasteelbook.jpg"
$ws.Range("C37").Value = "https://www.wog.ch/index.cfm/details/product/183350-Paprika-Steelbook-Edition-Blu-ray-UHD-2-Discs"
$ws.Range("D37").Value = "36 CHF"

$ws.Range("D38").Select()
